$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44162
$ws.Range("H2").Value = "Verde"
$ws.Range("J2").Value = 700
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1600
$ws.Range("P2").Value = 1600

# Row 3
$ws.Range("D3").Value2 = 44176
$ws.Range("J3").Value = 700

# Row 4
$ws.Range("D4").Value2 = 44179
$ws.Range("J4").Value = 200
$ws.Range("L4").Value = 1600
$ws.Range("M4").Value = 1600
$ws.Range("P4").Value = 1600

# Row 6
$ws.Range("D6").Value2 = 44473
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 2000

# Row 7
$ws.Range("D7").Value2 = 44168
$ws.Range("J7").Value = 200

# Row 8
$ws.Range("D8").Value2 = 44161
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 1700
$ws.Range("L8").Value = 1700
$ws.Range("M8").Value = 1700
$ws.Range("P8").Value = 1700

# Row 9
$ws.Range("D9").Value2 = 44165
$ws.Range("J9").Value = 300

# Row 10
$ws.Range("D10").Value2 = 44475
$ws.Range("J10").Value = 100

# Row 11
$ws.Range("D11").Value2 = 44474
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 2000
$ws.Range("P11").Value = 2000

# Row 12
$ws.Range("D12").Value2 = 44159
$ws.Range("J12").Value = 600
$ws.Range("L12").Value = 1700
$ws.Range("M12").Value = 1650
$ws.Range("P12").Value = 1650

# New row 13
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C13").Value = "Los Lagos"
$ws.Range("D13").Value2 = 44166
$ws.Range("D13").NumberFormat = $ws.Range("D12").NumberFormat
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Verde"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1600
$ws.Range("L13").Value = 1600
$ws.Range("M13").Value = 1600
$ws.Range("N13").Value = "$/kilo"
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1600
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
